# Auto-generated edit script: updates market-price derived columns (H-N)
# on the Leve profit sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 3203.7334
$ws.Range("I41").Value = 3033.2
$ws.Range("J41").Value = 3544.8
$ws.Range("K41").Value = 3033.2
$ws.Range("L41").Value = 3544.8
$ws.Range("M41").Value = -2593.2
$ws.Range("N41").Value = -4424.8
# Row 62
$ws.Range("H62").Value = 2203
$ws.Range("J62").Value = 1906
$ws.Range("L62").Value = 1906
$ws.Range("N62").Value = -3154
# Row 65
$ws.Range("H65").Value = 2203
$ws.Range("J65").Value = 1906
$ws.Range("L65").Value = 9530
$ws.Range("N65").Value = -15770
# Row 88
$ws.Range("H88").Value = 3283
$ws.Range("I88").Value = 2995
$ws.Range("J88").Value = 3331
$ws.Range("K88").Value = 2995
$ws.Range("L88").Value = 3331
$ws.Range("M88").Value = -2589
$ws.Range("N88").Value = -4143
# Row 91
$ws.Range("H91").Value = 3283
$ws.Range("I91").Value = 2995
$ws.Range("J91").Value = 3331
$ws.Range("K91").Value = 2995
$ws.Range("L91").Value = 3331
$ws.Range("M91").Value = -1591
$ws.Range("N91").Value = -6139
# Row 113
$ws.Range("H113").Value = 1999.875
$ws.Range("I113").Value = 1999.8572
$ws.Range("K113").Value = 1999.8572
$ws.Range("M113").Value = 1254.1428
# Row 116
$ws.Range("H116").Value = 30000
$ws.Range("I116").Value = 30000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 30000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -26558
$ws.Range("N116").ClearContents()
# Row 132
$ws.Range("H132").Value = 4270.1
$ws.Range("I132").Value = 4675.1113
$ws.Range("K132").Value = 14025.3339
$ws.Range("M132").Value = -11495.3339

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2004
$ws.Range("I45").Value = 2004
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2004
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1627
$ws.Range("N45").ClearContents()
# Row 122
$ws.Range("H122").Value = 3402.25
$ws.Range("J122").Value = 3256
$ws.Range("L122").Value = 9768
$ws.Range("N122").Value = -14668

$ws = $wb.Worksheets.Item("BSM")
# Row 123
$ws.Range("H123").Value = 93499.5
$ws.Range("J123").Value = 93499.5
$ws.Range("L123").Value = 93499.5
$ws.Range("N123").Value = -103299.5

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 322.25
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 444.5
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 444.5
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -1144.5
# Row 58
$ws.Range("H58").Value = 5766.1113
$ws.Range("I58").Value = 1298.3334
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 1298.3334
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = -1095.3334
$ws.Range("N58").Value = -8406
# Row 69
$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14251
# Row 72
$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -41256
# Row 136
$ws.Range("H136").Value = 5766.1113
$ws.Range("I136").Value = 1298.3334
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 3895.0002
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -1345.0002
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
# Row 49
$ws.Range("H49").Value = 250
$ws.Range("J49").Value = 250
$ws.Range("L49").Value = 750
$ws.Range("N49").Value = -1062
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# Row 113
$ws.Range("H113").Value = 414.66666
$ws.Range("J113").Value = 286
$ws.Range("L113").Value = 858
$ws.Range("N113").Value = -5198
# Row 141
$ws.Range("H141").Value = 3000
$ws.Range("J141").Value = 3000
$ws.Range("L141").Value = 9000
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 102
$ws.Range("H102").Value = 2979.8
$ws.Range("I102").Value = 2979.8
$ws.Range("K102").Value = 2979.8
$ws.Range("M102").Value = -1357.8
# Row 122
$ws.Range("H122").Value = 2177.5
$ws.Range("I122").Value = 2177.5
$ws.Range("K122").Value = 6532.5
$ws.Range("M122").Value = -4082.5
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 132
$ws.Range("H132").Value = 3812
$ws.Range("I132").Value = 3812
$ws.Range("K132").Value = 11436
$ws.Range("M132").Value = -8906

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4999
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 14997
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -19897
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# Row 136
$ws.Range("H136").Value = 19834.666
$ws.Range("I136").Value = 7252
$ws.Range("K136").Value = 21756
$ws.Range("M136").Value = -19206

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 3788.5
$ws.Range("I126").Value = 1718
$ws.Range("K126").Value = 5154
$ws.Range("M126").Value = -2684
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
